$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell's value to be stored as literal TEXT (mirrors the
# source data which is always an inline string), even when the text looks
# like a plain decimal number (e.g. "1.00", "27.60") that Excel would
# otherwise silently reinterpret/round as a Number on assignment.
# Sequence: flip to a text number-format so the incoming value is parsed as
# text, write the literal characters, then restore the "Normal" style so we
# don't leave a stray custom number-format applied to the cell.
function Set-TextValue($rng, [string]$val) {
    $rng.NumberFormat = "@"
    $rng.Characters().Text = $val
    $rng.Style = "Normal"
}

# ---- Column D (Price) updates ----
# Values that are NOT valid plain-decimal literals (contain thousands dots,
# subscript glyphs, etc.) are safe to assign directly - Excel keeps them as
# text automatically.
$ws.Range("D2").Value = "67.871.39"
$ws.Range("D3").Value = "2.627.80"
$ws.Range("D9").Value = "2.626.18"
$ws.Range("D16").Value = "3.104.40"
$ws.Range("D17").Value = "67.727.71"
$ws.Range("D18").Value = "2.630.86"
$ws.Range("D28").Value = "2.761.29"

# Values that DO look like plain decimal numbers need the text-forcing helper.
Set-TextValue $ws.Range("D5") "596.50"
Set-TextValue $ws.Range("D6") "153.76"
Set-TextValue $ws.Range("D10") "0.136"
Set-TextValue $ws.Range("D14") "27.60"
Set-TextValue $ws.Range("D20") "368.95"
Set-TextValue $ws.Range("D21") "7.43"
Set-TextValue $ws.Range("D25") "72.00"
Set-TextValue $ws.Range("D27") "9.84"
Set-TextValue $ws.Range("D30") "1.00"
Set-TextValue $ws.Range("D31") "574.54"
Set-TextValue $ws.Range("D35") "0.999"
Set-TextValue $ws.Range("D37") "1.54"
Set-TextValue $ws.Range("D38") "158.68"
Set-TextValue $ws.Range("D39") "19.15"
Set-TextValue $ws.Range("D47") "40.03"
Set-TextValue $ws.Range("D48") "155.36"
Set-TextValue $ws.Range("D50") "21.97"

# ---- Column E (Volume(1h)) updates ----
# All values share the "  +x.xx%  " padded layout, which is never parsed as
# a numeric percentage by Excel, so plain assignment keeps them as text.
$ws.Range("E2").Value  = "  +0.97%  "
$ws.Range("E3").Value  = "  +0.53%  "
$ws.Range("E4").Value  = "  -0.04%  "
$ws.Range("E5").Value  = "  +0.30%  "
$ws.Range("E6").Value  = "  +0.92%  "
$ws.Range("E7").Value  = "  +0.00%  "
$ws.Range("E9").Value  = "  +0.52%  "
$ws.Range("E10").Value = "  +10.52%  "
$ws.Range("E11").Value = "  -0.55%  "
$ws.Range("E12").Value = "  +1.20%  "
$ws.Range("E13").Value = "  +0.80%  "
$ws.Range("E14").Value = "  +0.57%  "
$ws.Range("E15").Value = "  +4.91%  "
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("E17").Value = "  +0.92%  "
$ws.Range("E18").Value = "  +1.20%  "
$ws.Range("E19").Value = "  +3.74%  "
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("E21").Value = "  +0.86%  "
$ws.Range("E22").Value = "  -1.30%  "
$ws.Range("E23").Value = "  -1.21%  "
$ws.Range("E24").Value = "  +1.76%  "
$ws.Range("E25").Value = "  +8.33%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  -1.53%  "
$ws.Range("E28").Value = "  +0.87%  "
$ws.Range("E29").Value = "  +3.30%  "
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("E31").Value = "  -0.78%  "
$ws.Range("E32").Value = "  +2.31%  "
$ws.Range("E33").Value = "  +1.40%  "
$ws.Range("E34").Value = "  +1.41%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  +5.19%  "
$ws.Range("E37").Value = "  +2.28%  "
$ws.Range("E39").Value = "  +0.79%  "
$ws.Range("E40").Value = "  +4.77%  "
$ws.Range("E41").Value = "  +0.49%  "
$ws.Range("E42").Value = "  +2.34%  "
$ws.Range("E45").Value = "  +6.13%  "
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("E47").Value = "  -2.48%  "
$ws.Range("E48").Value = "  -0.36%  "
$ws.Range("E49").Value = "  -0.75%  "
$ws.Range("E50").Value = "  +1.38%  "
$ws.Range("E51").Value = "  +0.05%  "

# ---- Rows 43 / 44: dogwifhat and BabyDogeCoin swap rank order ----
# Row 43 becomes BabyDogeCoin, row 44 becomes dogwifhat.
$ws.Range("B43").Value = "BabyDogeCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D43").Value = "0.0₆0334"
$ws.Range("E43").Value = "  +16.57%  "

$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D44") "2.63"
$ws.Range("E44").Value = "  +3.90%  "
